# Update Issue 357 Merge and resolve conflicts
#
# 1) Two diagram rectangles that were labelled "CoordData" are renamed to
#    "InstructorData" (slide 1, shapes #2 and #35 in z-order).
# 2) The cached "datetimeFigureOut" date placeholder text ("7/8/2012") on the
#    slide master and on every one of its 11 slide layouts is refreshed to
#    "12/1/2012" (this is what PowerPoint does to the auto date field's
#    cached display text on save).

$p = $ppt.ActivePresentation

# --- 1) Rename the "CoordData" shapes to "InstructorData" --------------
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq "CoordData") {
            $shp.TextFrame.TextRange.Text = "InstructorData"
        }
    }
}

# --- 2) Refresh the cached date placeholder text ------------------------
$oldDate = "7/8/2012"
$newDate = "12/1/2012"

$m = $p.SlideMaster

function Update-DatePlaceholder($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

Update-DatePlaceholder $m.Shapes

for ($k = 1; $k -le $m.CustomLayouts.Count; $k++) {
    $layout = $m.CustomLayouts.Item($k)
    Update-DatePlaceholder $layout.Shapes
}
